# "Add files via upload" — turns the single "Room Management" testcase
# sheet into two sprint sections ("Sprint 1" / "Sprint 2"), each with its
# own centered title bar, header block and testcase table. Sprint 2 is a
# new set of testcases (new author "Simon Leu", new date, new rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New row 1: "Sprint 1" title, merged B1:D1, centered.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Sprint 1"
$ws.Range("B1:D1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B1:D1").Merge()

# ---------------------------------------------------------------------
# 2. Rows 8-10: row height only changes (content identical). Recomputed
#    by Excel after the default font/row height shifted; reproduce the
#    same 0.96 scale it ended up with.
# ---------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 28.8
$ws.Rows.Item(10).RowHeight = 43.2

# ---------------------------------------------------------------------
# 3. Clear the leftover placeholder formatting on what becomes the
#    blank separator row between the two sprint blocks (old row 17 had
#    styled-but-empty C/D cells; in the final sheet it is fully blank).
# ---------------------------------------------------------------------
$ws.Range("C17:D17").Clear()

# ---------------------------------------------------------------------
# 4. Row 13: "Sprint 2" title, merged B13:D13, centered (same recipe as
#    row 1).
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "Sprint 2"
$ws.Range("B13:D13").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B13:D13").Merge()

# ---------------------------------------------------------------------
# 5. Row 14: "Testcases " title (same shaded style as B2). Use
#    Copy/PasteSpecial(formats) to clone the exact style — direct
#    `.Style = other.Style` assignment doesn't transfer formatting in
#    this host.
# ---------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B14").Value = "Testcases "
$ws.Range("C14:D14").Clear()   # drop the leftover placeholder styling

# ---------------------------------------------------------------------
# 6. Row 15: project / version header line (same layout/style as row 3).
# ---------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "Projekt:"

$ws.Range("C15").Value = "Room Management"

$ws.Range("D3").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = "Version:"

# "1.0" looks numeric, so a plain .Value assignment would coerce it to
# the number 1 — paste the existing "1.0" text cell's value instead so
# it stays a shared-string text cell (matches the source row 3 E-cell).
$ws.Range("E3").Copy()
$ws.Range("E15").PasteSpecial(-4163)   # xlPasteValues

# ---------------------------------------------------------------------
# 7. Row 16: creation-date / author header line (same layout/style as
#    row 4), with the new date and the new author.
# ---------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = "Erstellungsdatum:"

$ws.Range("C4").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 43286

$ws.Range("D4").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Ersteller:"

$ws.Range("E16").Value = "Simon Leu"

# ---------------------------------------------------------------------
# 8. Row 18: table header (Testfallnummer / Beschreibung / Erwartetes
#    Ergebnis), same style as row 6.
# ---------------------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "Testfallnummer "

$ws.Range("C6").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "Beschreibung"

$ws.Range("D6").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Erwartetes Ergebnis"

# ---------------------------------------------------------------------
# 9. Rows 19-22: the four Sprint-2 testcases. Numbering column is
#    vertical-centered; description/result columns are vertical-centered
#    + wrap (except C22, which only got vertical-center in the source
#    file).
# ---------------------------------------------------------------------
$ws.Range("B19").Value = 1
$ws.Range("B19").VerticalAlignment = -4108   # xlCenter

$ws.Range("C19").Value = "Der Benutzer ruft die Raumliste auf."
$ws.Range("C19").VerticalAlignment = -4108
$ws.Range("C19").WrapText = $true

$ws.Range("D19").Value = "Es werden alle verfügbaren Räume angezeigt."
$ws.Range("D19").VerticalAlignment = -4108
$ws.Range("D19").WrapText = $true

$ws.Range("B20").Value = 2
$ws.Range("B20").VerticalAlignment = -4108

$ws.Range("C20").Value = "Der Benutzer ruft die Detail eines Raumes auf."
$ws.Range("C20").VerticalAlignment = -4108
$ws.Range("C20").WrapText = $true

$ws.Range("D20").Value = "Es wird eine Detailansicht mit der Verfügbarkeit des Raumes angezeigt."
$ws.Range("D20").VerticalAlignment = -4108
$ws.Range("D20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 28.8

$ws.Range("B21").Value = 3
$ws.Range("B21").VerticalAlignment = -4108

$ws.Range("C21").Value = "Der Benutzer sucht auf der Raumliste die verfügbaren Räume für einen Tag und eine Zeit."
$ws.Range("C21").VerticalAlignment = -4108
$ws.Range("C21").WrapText = $true

$ws.Range("D21").Value = "Alle zu dieser Zeit an diesem Tag verfügbaren Räume werden angezeigt."
$ws.Range("D21").VerticalAlignment = -4108
$ws.Range("D21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 28.8

$ws.Range("B22").Value = 4
$ws.Range("B22").VerticalAlignment = -4108

$ws.Range("C22").Value = "Der Benutzer reserviert einen Raum an einem freien Datum."
$ws.Range("C22").VerticalAlignment = -4108

$ws.Range("D22").Value = "Der Benutzer erhält eine Bestätigung seiner Reservierung und es wird eine Reservation erstellt."
$ws.Range("D22").VerticalAlignment = -4108
$ws.Range("D22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 43.2

# ---------------------------------------------------------------------
# 10. View state: scroll so the new bottom rows are visible, select the
#     last-edited cell (matches the saved selection in the source file).
# ---------------------------------------------------------------------
$ws.Range("D22").Select()
